# Insert a new row at row 7 (pushes existing rows 7..100 down to 8..101)
# and populate it with the new data point, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44685
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112021
$ws.Range("G7").Value = "Ají"
$ws.Range("H7").Value = "Inferno"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 26000
$ws.Range("L7").Value = 27000
$ws.Range("M7").Value = 26600
$ws.Range("N7").Value = "$/caja 15 kilos"
$ws.Range("O7").Value = "Provincia de Huasco"
$ws.Range("P7").Value = 1773
$ws.Range("Q7").Value = 15
$ws.Range("R7").Value = "Hortaliza"
